$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Maggie Burton got +1 hour on her timecard: 19h 30m -> 20h 30m
$ws.Range("B7").Value = "20h 30m"

# Reflect the final cell selection left in the saved file
$ws.Range("D11").Select()
